# Daily scrape update - 2026-01-30 04:13:01 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = "'1331509"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1331509"
$ws.Range("C2").Value = "Commercial area internship - KAM Jr."
$ws.Range("D2").Value = "Antiguo Cuscatlán, El Salvador"
$ws.Range("E2").Value = "No"
$ws.Range("F2").Value = "0 applicants"
$ws.Range("G2").Value = "3 - 6 Months"
$ws.Range("H2").Value = "Samsung Electronics Latinoamerica"

# --- Row 3 ---
$ws.Range("A3").Value = "'1325830"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1325830"
$ws.Range("C3").Value = "Marketing Intern"
$ws.Range("D3").Value = "Makati City, Metro Manila, Philippines"
$ws.Range("E3").Value = "No"
$ws.Range("F3").Value = "42 applicants"
$ws.Range("G3").Value = "3 - 6 Months"
$ws.Range("H3").Value = "Consistent Frozen Solutions Corporation"

# --- Row 4 ---
$ws.Range("A4").Value = "'1322997"
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1322997"
$ws.Range("C4").Value = "Marketing"
$ws.Range("D4").Value = "Sfax, Tunisie"
$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "8 applicants"
$ws.Range("G4").Value = "9 - 12 Weeks"
$ws.Range("H4").Value = "English Cultural Center"

# --- Row 5 ---
$ws.Range("A5").Value = "'1316723"
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1316723"
$ws.Range("C5").Value = "[Impact Porto Alegre]- Social Media"
$ws.Range("D5").Value = "Porto Alegre, RS, Brasil"
$ws.Range("E5").Value = "No"
$ws.Range("F5").Value = "84 applicants"
$ws.Range("G5").Value = "9 - 12 Weeks"
$ws.Range("H5").Value = "ESCOLA GIORDANO BRUNO LTDA"

# --- Row 6 ---
$ws.Range("A6").Value = "'1299853"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = "https://aiesec.org/opportunity/global-talent/1299853"
$ws.Range("C6").Value = "[Impact Porto Alegre]- Social Media"
$ws.Range("D6").Value = "Bom Fim, Brazil"
$ws.Range("E6").Value = "No"
$ws.Range("F6").Value = "147 applicants"
$ws.Range("G6").Value = "9 - 12 Weeks"
$ws.Range("H6").Value = "ESCOLA GIORDANO BRUNO LTDA"

# --- Remove old rows 7-10 (data no longer present in the latest scrape) ---
$ws.Range("A7:H10").EntireRow.Delete()

# --- Column width adjustments ---
# (ColumnWidth is specified in Excel "characters"; Excel stores widths with a
#  constant +5/6 padding offset for the default font, so we subtract that
#  offset here to land on the exact target widths in the saved file.)
$offset = 0.8333333333333333
$ws.Columns.Item(3).ColumnWidth = 39 - $offset
$ws.Columns.Item(4).ColumnWidth = 41 - $offset
$ws.Columns.Item(7).ColumnWidth = 15 - $offset
$ws.Columns.Item(8).ColumnWidth = 42 - $offset
